# "Generate Report for Handback"
# Stamps fresh handoff/handback timestamps into the status workbook as each
# locale file is (re-)processed, mirroring the real report generator:
#   1) Overview sheet: refresh "Latest HO Xliff Generate Date" for the
#      cb4dc240 markdown file.
#   2) zh-cn sheet: refresh the Correspond Handoff / Handback datetimes for
#      the cb4dc240 xliff row.
#   3) de-de sheet: refresh the Correspond Handoff datetime for BOTH rows
#      (they pick up the same "now" stamp used for the Overview refresh),
#      plus the Correspond Handback datetime for the cb4dc240 row.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$now1 = "2016-08-30 06:51:16"

# 1) Overview!G3 - Latest HO Xliff Generate Date (cb4dc240 row)
$overview.Cells.Item(3, 7).Value = $now1

# 2) zh-cn!H3 / K3 - Correspond Handoff / Handback Datetime (cb4dc240 row)
$zhcn.Cells.Item(3, 8).Value = "2016-08-30 06:51:11"
$zhcn.Cells.Item(3, 11).Value = "2016-08-30 06:51:28"

# 3) de-de!H2 / H3 - Correspond Handoff Datetime (both rows, same stamp as Overview)
$dede.Cells.Item(2, 8).Value = $now1
$dede.Cells.Item(3, 8).Value = $now1

# de-de!K3 - Correspond Handback Datetime (cb4dc240 row)
$dede.Cells.Item(3, 11).Value = "2016-08-30 06:51:35"
